$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (toy_text_practice_1 / question 1): replace the question text
# with a new, deliberately long, test question.
$ws.Range("F5").Value = "This is a very long question that is completely useless, but we need it to test how a two-line question looks like!"

# --- Row 3 (toy_text_2 / question 1): fix a typo in the question text
# ("infred" -> "infered") ...
$ws.Range("F3").Value = "What can be infered from eye-movements?"

# ... and add the missing distractor columns (they were empty before).
$ws.Range("I3").Value = "whatever"
$ws.Range("J3").Value = "watever"
$ws.Range("K3").Value = "and this is another very long answer option to test what happens"

# --- Row 4 (toy_text_2 / question 2): add the missing distractor_3 value,
# matching the wrap-text look of its neighbouring span-text cells.
$ws.Range("K4").Value = "whatever"
$ws.Range("K4").WrapText = $true

# Re-apply wrap formatting on L4/M4 so they pick up the same (de-duplicated)
# style as K4 instead of their old, redundant bordered style.
$ws.Range("L4").WrapText = $true
$ws.Range("M4").WrapText = $true

# Leave the cursor where the author ended up after making these edits.
$ws.Range("H13").Select() | Out-Null
